# timing of mask + face presentation
# Sheet1 held timing/response data that got re-collected; update the
# A2:B5 block with the refreshed values (B2 now blank -> no response
# recorded for that trial) and keep the header row as-is.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A2").Value = 1
$ws1.Range("B2").ClearContents()

$ws1.Range("A3").Value = 0
$ws1.Range("B3").Value = 4

$ws1.Range("A4").Value = 0
$ws1.Range("B4").Value = 9

$ws1.Range("A5").Value = 1
$ws1.Range("B5").Value = 9

$ws1.Columns.Item(1).ColumnWidth = 15.5703125
$ws1.Columns.Item(2).ColumnWidth = 15.5703125
